$d = $word.ActiveDocument

# --- 1. Remove the existing "_GoBack" bookmark from the last paragraph.
#        It will be re-created at the end of the newly appended content. ---
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
} catch {
    # no pre-existing _GoBack bookmark - nothing to remove
}

# --- 2. Build the OOXML fragment for the new paragraphs that must be
#        appended after "...viene spento automaticamente." ---
$newParagraphsXml =
    '<w:p><w:pPr><w:jc w:val="both"/></w:pPr></w:p>' +
    '<w:p><w:pPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Ridurre i tempi (30s </w:t></w:r>' +
    '<w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r>' +
    '<w:r><w:t xml:space="preserve"> 5s e 20m</w:t></w:r>' +
    '<w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r>' +
    '<w:r><w:t>1m) per permettere il verificarsi di tutti gli eventi possibili durante la discussione.</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Variare di temperatura/umidit&#224; va fatto tramite codice (usando random()) perch&#233; </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>cooja</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> non lo fa.</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:t>Potrebbe prevedere una variabile di stato lato CU che, una volta mostrata la lista dei comandi, permetta all&#8217;utente di conoscere se il nodo 4 &#232; impostato a sauna/turco. (La CU conosce questo stato perch&#233; le viene comunicato da nodo 4).</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:jc w:val="both"/></w:pPr></w:p>' +
    '<w:p><w:pPr><w:jc w:val="both"/></w:pPr></w:p>' +
    '<w:p><w:pPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Consiglio: mettere il controllo </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>alarm</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">=1 sulla CU; mettere </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>etimer</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> a 2 secondi.</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'

$packageXml =
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
    $newParagraphsXml +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- 3. Insert the fragment at the very end of the document body. ---
$insertionPoint = $d.Range($d.Content.End, $d.Content.End)
$insertionPoint.InsertXML($packageXml)
